# Update the "Base" sheet's regression rows so the hard-coded Salesforce
# CreatedById ('005q0000003GGfP') used to look up freshly-created records is
# replaced with the new user id '0051I000000guzk'. Four cells on the Base
# sheet embed that id: the fetchwo (B2), fetchaccount (E2) and
# FetchServiceContract (G2) SOQL strings, plus the raw Username value (J2).

$wb = $excel.ActiveWorkbook
$wsBase = $wb.Worksheets.Item("Base")
$wsDev = $wb.Worksheets.Item("Dev")

$oldId = "005q0000003GGfP"
$newId = "0051I000000guzk"

$wsBase.Range("B2").Value = "Select Name , Id,FORMAT(SVMXC__Actual_Initial_Response__c), FORMAT(SVMXC__Actual_Onsite_Response__c) from SVMXC__Service_Order__c where Createdby.Id = '" + $newId + "' Order by CreatedDate DESC Limit 1"
$wsBase.Range("E2").Value = "Select Name , Id from Account where Name = 'Account Automation 1' AND Createdby.Id = '" + $newId + "' Order by CreatedDate DESC Limit 1"
$wsBase.Range("G2").Value = "Select Name , Id from SVMXC__Service_Contract__c where Createdby.Id = '" + $newId + "' Order by CreatedDate DESC Limit 1"
# Leading apostrophe keeps this a text entry (and keeps the existing
# quote-prefix cell style) instead of Excel re-guessing the format.
$wsBase.Range("J2").Value = "'" + $newId

# The regression sheet for Dev grew a couple of extra wrapped lines, so bump
# its row-2 height to keep the longer cell text fully visible.
$wsDev.Rows.Item(2).RowHeight = 120

# Make "Base" the active sheet/tab again, with the cursor parked on B7.
$wsBase.Activate()
$wsBase.Range("B7").Select()
